$d = $word.ActiveDocument

# The document ends with an empty trailing paragraph just before the
# sectPr. Turn it into the "Requirement Lead sign off:" paragraph, then
# add a new paragraph after it for "Vincent Lam" with a first-line indent.
$signOffPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$signOffPara.Range.Text = "Requirement Lead"

$insertPos = $signOffPara.Range.End - 1
$tailRange = $d.Range($insertPos, $insertPos)
$tailRange.InsertAfter(" sign off:")

$signOffPara.Range.InsertParagraphAfter()

$namePara = $d.Paragraphs.Item($d.Paragraphs.Count)
$namePara.Range.Text = "Vincent Lam"
$namePara.FirstLineIndent = 36

# Explicitly (re)assert the page as portrait oriented.
$d.PageSetup.Orientation = 0
